$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Students")

$ws.Range("A7").Value = "R-22-0124"
$ws.Range("B7").Value = "r220124@famt.ac.in"
$ws.Range("C7").Value = "Mansi Surendra Agre"
$ws.Range("D7").Value = "abcd@gmail.com"
$ws.Range("E7").Value = '$2b$12$NH5vRBM2YWaelWJPnz.K4.gj1i5TQo3qXdUURWAKRekdm9uGxumSO'
$ws.Range("F7").Value = "student"
